# Regenerate the "K" column (column G) values using updated calculations
# (K = computed from strike count data instead of the old "Strike#" based value).
# This mirrors the upstream commit that recalculated std/mean and rewrote the
# K column (G2:G26) with freshly computed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 4
    3  = 2
    4  = 6
    5  = 3
    6  = 1
    7  = 3
    8  = 6
    9  = 3
    10 = 4
    11 = 6
    12 = 6
    13 = 2
    14 = 3
    15 = 6
    16 = 2
    17 = 7
    18 = 4
    19 = 5
    20 = 4
    21 = 5
    22 = 2
    23 = 3
    24 = 4
    25 = 5
    26 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
